$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "https://sites.research.unimelb.edu.au/research-funding/domestic", "Accepting and managing your funding"),
    @(3, "https://sites.research.unimelb.edu.au/research-funding/domestic/manage-your-grant", "Accepting, and, funding, managing, your"),
    @(4, "https://sites.research.unimelb.edu.au/research-funding/domestic/cancer-council-victoria-grants-in-aid", "Accepting, and, funding, managing, your"),
    @(5, "https://sites.research.unimelb.edu.au/research-funding/domestic/vca-grants", "Accepting, and, funding, managing, your"),
    @(6, "https://sites.research.unimelb.edu.au/research-funding/domestic/victorian-cancer-agency-grants2", "Accepting, and, funding, managing, your"),
    @(7, "https://sites.research.unimelb.edu.au/research-funding/domestic/schemes-for-international-collaboration", "Accepting, and, funding, managing, your"),
    @(8, "https://sites.research.unimelb.edu.au/research-funding/domestic/westpac", "Accepting, and, funding, managing, your"),
    @(9, "https://sites.research.unimelb.edu.au/research-funding/domestic/djpr-victorian-medical-research-acceleration-fund-vmraf-and-mrna-vmraf", "Accepting, and, funding, managing, your"),
    @(10, "https://sites.research.unimelb.edu.au/research-funding/domestic/mrna-vraf", "Accepting, and, funding, managing, your")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
